# "map updates changed the order of almost all the maps in the game
#  updated memory requirements of maps"
#
# The "maps" block used to live in rows 121-136 of Sheet1. Map entries
# (index in col A, memory in col B, comments in col C/D/F/H/I/K/L) get
# reordered, map #1's description/memory is updated, and the old map #2-5
# entries (which used to sit right after #1) move down to the very end of
# the block (now rows 137-140), growing the block from 16 to 20 rows
# (121-140 instead of 121-136).
#
# Simplest faithful reproduction: wipe the old block and rewrite it with
# the final values, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old "maps" block (rows 121-136) entirely so no stale cells
# from the previous layout/ordering linger around.
$ws.Rows("121:136").Delete()

# Header row (unchanged content, just re-written since the row was wiped).
$ws.Range("B121").Value = "normal"
$ws.Range("C121").Value = "subs"
$ws.Range("D121").Value = "reprogram"

# Map #7
$ws.Range("A126").Value = 7
$ws.Range("B126").Value = 64
$ws.Range("F126").Value = "very easy, should be an earlier level"

# Map #6
$ws.Range("A127").Value = 6
$ws.Range("B127").Value = 78
$ws.Range("F127").Value = "very easy, should be an earlier level"

# Map #10
$ws.Range("A128").Value = 10
$ws.Range("B128").Value = 154
$ws.Range("F128").Value = "easy - very straightforward"

# Map #8 (now also carries the "jump-move forward..." / "sub1 = ..." notes)
$ws.Range("A129").Value = 8
$ws.Range("B129").Value = 156
$ws.Range("F129").Value = "fairly easy - fun water map"
$ws.Range("H129").Value = "jump-move forward and jump-move forward 3 are EASILY mass repeated on this map, will test with subs"
$ws.Range("I129").Value = "sub1 = jump/move forward, sub2 = jump"

# Map #9 (memory requirement bumped, note moved off to map #8 above)
$ws.Range("A130").Value = 9
$ws.Range("B130").Value = 306
$ws.Range("C130").Value = 164
$ws.Range("F130").Value = "fun - lots of jumping"

# Map #1 - updated memory requirement and description
$ws.Range("A131").Value = 1
$ws.Range("B131").Value = 56
$ws.Range("F131").Value = "short and sweet"

# Map #11
$ws.Range("A132").Value = 11
$ws.Range("B132").Value = 128
$ws.Range("F132").Value = "not very complex, fairly easy as well - lots of random extra stuff not related to finishing the map"

# Map #12 (now also flagged "can't be beaten")
$ws.Range("A133").Value = 12
$ws.Range("F133").Value = "can't be beaten"

# Map #13
$ws.Range("A134").Value = 13
$ws.Range("B134").Value = 182
$ws.Range("F134").Value = "could probably lose the reprogram square, additionally, not sure if intended, but the last two switches can be skipped entirely, may want to disable jump on this map, or make some reason for the switches"

# Map #14
$ws.Range("A135").Value = 14
$ws.Range("B135").Value = 318
$ws.Range("F135").Value = "very linear, interesting figuring out what does what"
$ws.Range("L135").Value = "needs edge squares removed maybe"

# Map #15
$ws.Range("A136").Value = 15
$ws.Range("F136").Value = "can't be beaten"

# Map #2 (moved down from the top of the block)
$ws.Range("A137").Value = 2
$ws.Range("B137").Value = 156
$ws.Range("F137").Value = "annoying as hell"

# Map #3
$ws.Range("A138").Value = 3
$ws.Range("B138").Value = 110
$ws.Range("D138").Value = "62 with reprogram"
$ws.Range("F138").Value = "annoying as hell"

# Map #4
$ws.Range("A139").Value = 4
$ws.Range("B139").Value = 88
$ws.Range("F139").Value = "annoying as hell"
$ws.Range("K139").Value = "used left switch"

# Map #5
$ws.Range("A140").Value = 5
$ws.Range("B140").Value = 92
$ws.Range("F140").Value = "annoying as hell"
$ws.Range("K140").Value = "fixed bad switch"

# Selection follows the edit (matches the final cursor position recorded
# in the workbook).
[void]$ws.Range("F132").Select()
